# Add "2022-Q3" sheet (fund holdings) right after "总计", shifting all
# quarterly sheets one slot to the right, and add a matching summary row
# at the top of the "总计" sheet's data.

$wb = $excel.ActiveWorkbook
$wsTotal = $wb.Worksheets.Item("总计")
$wsQ2 = $wb.Worksheets.Item("2022-Q2")

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q3" sheet by copying "2022-Q2" (so it inherits
#    the same column layout/styles), placed right after "总计".
# ---------------------------------------------------------------------
$wsQ2.Copy($null, $wsTotal)
$wsQ3 = $wb.Worksheets.Item(2)
$wsQ3.Name = "2022-Q3"

# The copied sheet only has 2 data rows (rows 2-3); the new quarter needs
# 3, so give row 4 the same formatting as row 2 before filling it in.
$wsQ3.Cells.Item(2,1).Copy()
$wsQ3.Cells.Item(4,1).PasteSpecial(-4122)

# Columns B, D, E, F, G hold text (fund codes / formatted numbers as
# strings) in every quarterly sheet - force Text format so Excel doesn't
# auto-convert them to numbers (dropping leading zeros / reformatting).
$wsQ3.Range("B2:B4").NumberFormat = "@"
$wsQ3.Range("D2:G4").NumberFormat = "@"

$wsQ3.Cells.Item(2,1).Value = 0
$wsQ3.Cells.Item(2,2).Value = "011355"
$wsQ3.Cells.Item(2,3).Value = "华泰柏瑞港股通时代机遇混合A"
$wsQ3.Cells.Item(2,4).Value = "0.54"
$wsQ3.Cells.Item(2,5).Value = "91.80"
$wsQ3.Cells.Item(2,6).Value = "4.59"
$wsQ3.Cells.Item(2,7).Value = "0.0248"
$wsQ3.Cells.Item(2,8).Value = 9

$wsQ3.Cells.Item(3,1).Value = 1
$wsQ3.Cells.Item(3,2).Value = "460010"
$wsQ3.Cells.Item(3,3).Value = "华泰柏瑞亚洲领导企业混合（QDII）"
$wsQ3.Cells.Item(3,4).Value = "0.36"
$wsQ3.Cells.Item(3,5).Value = "93.44"
$wsQ3.Cells.Item(3,6).Value = "4.90"
$wsQ3.Cells.Item(3,7).Value = "0.0176"
$wsQ3.Cells.Item(3,8).Value = 10

$wsQ3.Cells.Item(4,1).Value = 2
$wsQ3.Cells.Item(4,2).Value = "011356"
$wsQ3.Cells.Item(4,3).Value = "华泰柏瑞港股通时代机遇混合C"
$wsQ3.Cells.Item(4,4).Value = "0.24"
$wsQ3.Cells.Item(4,5).Value = "91.80"
$wsQ3.Cells.Item(4,6).Value = "4.59"
$wsQ3.Cells.Item(4,7).Value = "0.0110"
$wsQ3.Cells.Item(4,8).Value = 9

# ---------------------------------------------------------------------
# 2) Insert a new summary row for 2022-Q3 at the top of the "总计" sheet.
# ---------------------------------------------------------------------
$wsTotal.Rows.Item(2).Insert()
$wsTotal.Range("B2:D2").ClearFormats()

# Row 2 lost column A's style when it was freshly inserted - copy it from
# row 3 (the shifted-down former row 2) which still carries it.
$wsTotal.Cells.Item(3,1).Copy()
$wsTotal.Cells.Item(2,1).PasteSpecial(-4122)

$wsTotal.Cells.Item(2,1).Value = 0
$wsTotal.Cells.Item(2,2).Value = "2022-Q3"
$wsTotal.Cells.Item(2,3).Value = 3
$wsTotal.Cells.Item(2,4).Value = 0.05

# The other rows kept their old data (shifted down 1) - just renumber the
# sequential index column A.
for ($r = 3; $r -le 8; $r++) {
    $wsTotal.Cells.Item($r,1).Value = $r - 2
}
